# Second trial for MQTT #71
#
# 1) Re-cache the "datetimeFigureOut" date placeholder (slide master + all
#    slide layouts) from 31/12/2022 to 01/01/2023.
# 2) Extend the "Rectangle 11" description text box on slide 1 to mention
#    that the communication now happens via (bold) MQTT.

$p = $ppt.ActivePresentation
$newDate = "01/01/2023"

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# -- slide master --------------------------------------------------------
$master = $p.SlideMaster
Update-DatePlaceholder $master

# -- every slide layout under the master ---------------------------------
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout
}

# -- slide 1: "Rectangle 11" description textbox -------------------------
$slide = $p.Slides.Item(1)
$descShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 11") {
        $descShape = $shp
    }
}

$tr = $descShape.TextFrame.TextRange

# The last two runs of the paragraph currently read:
#   ", but integrates the Python service code into a synchronous service
#    passing its data on to a synchronous sink. The remaining properties
#    of the example are " + "the same."
# Replace that tail (164 characters, starting right after "examples.python")
# with the extended sentence.
$oldTail = ", but integrates the Python service code into a synchronous service passing its data on to a synchronous sink. The remaining properties of the example are the same."
$tailStart = $tr.Text.IndexOf($oldTail) + 1
$tailLen = $oldTail.Length

$newTail = ", but integrates the Python service code into a synchronous service passing its data on to a synchronous sink. The remaining properties of the example are the same except for that the communication happens with MQTT."
$tr.Characters($tailStart, $tailLen).Text = $newTail

# Bold just "MQTT" within the newly-inserted text.
$fullText = $descShape.TextFrame.TextRange.Text
$mqttStart = $fullText.IndexOf("MQTT") + 1
$tr.Characters($mqttStart, 4).Font.Bold = $true

Write-Host "Updated date placeholders and MQTT description."
